$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.277.45"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.359.90"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.52"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.74"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("E10").Value = "  +4.87%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.45"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").Value = "2.781.59"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "57.287.67"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "2.369.85"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.57"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "330.35"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.44"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.97"
$ws.Range("E24").Value = "  +14.82%  "
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +11.40%  "
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.65"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  -4.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.03"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  +5.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.74"
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.81"
$ws.Range("E40").Value = "  +7.31%  "
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "290.90"
$ws.Range("E42").Value = "  +4.19%  "
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0939"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.565"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0219"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.16"
$ws.Range("E49").Value = "  +4.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.68"
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("B51").Value = "Polygon"
$ws.Range("C51").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.362"
$ws.Range("E51").Value = "  -5.54%  "
